$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 367.61905
$ws.Range("I2").Value = 196.42857
$ws.Range("J2").Value = 453.2143
$ws.Range("K2").Value = 196.42857
$ws.Range("L2").Value = 453.2143
$ws.Range("M2").Value = -83.42857000000001
$ws.Range("N2").Value = -679.2143
$ws.Range("H8").Value = 1450
$ws.Range("I8").Value = 416.66666
$ws.Range("J8").Value = 3000
$ws.Range("K8").Value = 1249.99998
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = -1110.99998
$ws.Range("N8").Value = -9278
$ws.Range("H9").Value = 12910
$ws.Range("I9").Value = 15052.667
$ws.Range("J9").Value = 54
$ws.Range("K9").Value = 15052.667
$ws.Range("L9").Value = 54
$ws.Range("M9").Value = -14883.667
$ws.Range("N9").Value = -392
$ws.Range("H38").Value = 1325.2222
$ws.Range("I38").Value = 796.5
$ws.Range("K38").Value = 2389.5
$ws.Range("M38").Value = -2017.5
$ws.Range("H44").Value = 14025
$ws.Range("J44").Value = 14025
$ws.Range("L44").Value = 14025
$ws.Range("N44").Value = -14949
$ws.Range("H51").Value = 2814.913
$ws.Range("J51").Value = 3076.8
$ws.Range("L51").Value = 3076.8
$ws.Range("N51").Value = -4044.8
$ws.Range("H53").Value = 678.4545000000001
$ws.Range("I53").Value = 216.28572
$ws.Range("J53").Value = 1487.25
$ws.Range("K53").Value = 216.28572
$ws.Range("L53").Value = 1487.25
$ws.Range("M53").Value = 420.71428
$ws.Range("N53").Value = -2761.25
$ws.Range("H54").Value = 19997.666
$ws.Range("I54").Value = 19997.666
$ws.Range("K54").Value = 19997.666
$ws.Range("M54").Value = -19511.666
$ws.Range("H55").Value = 398.44446
$ws.Range("I55").Value = 364.83334
$ws.Range("J55").Value = 465.66666
$ws.Range("K55").Value = 364.83334
$ws.Range("L55").Value = 465.66666
$ws.Range("M55").Value = -150.83334
$ws.Range("N55").Value = -893.66666
$ws.Range("H59").Value = 524.5
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H61").Value = 102.6
$ws.Range("I61").Value = 102.6
$ws.Range("K61").Value = 307.8
$ws.Range("M61").Value = -135.8
$ws.Range("H87").Value = 88535.664
$ws.Range("J87").Value = 88535.664
$ws.Range("L87").Value = 88535.664
$ws.Range("N87").Value = -91031.664
$ws.Range("H90").Value = 88535.664
$ws.Range("J90").Value = 88535.664
$ws.Range("L90").Value = 265606.992
$ws.Range("N90").Value = -278086.992
$ws.Range("H112").Value = 40068.57
$ws.Range("J112").Value = 41497.035
$ws.Range("L112").Value = 124491.105
$ws.Range("N112").Value = -126707.105
$ws.Range("H138").Value = 7695348
$ws.Range("I138").Value = 1639.3334
$ws.Range("J138").Value = 9437320
$ws.Range("K138").Value = 4918.0002
$ws.Range("L138").Value = 28311960
$ws.Range("M138").Value = 221.9997999999996
$ws.Range("N138").Value = -28322240

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 506
$ws.Range("I2").Value = 356.88235
$ws.Range("J2").Value = 1013
$ws.Range("K2").Value = 356.88235
$ws.Range("L2").Value = 1013
$ws.Range("M2").Value = -243.88235
$ws.Range("N2").Value = -1239
$ws.Range("H74").Value = 6037.5117
$ws.Range("I74").Value = 1241.5
$ws.Range("K74").Value = 1241.5
$ws.Range("M74").Value = -367.5
$ws.Range("H77").Value = 6037.5117
$ws.Range("I77").Value = 1241.5
$ws.Range("K77").Value = 6207.5
$ws.Range("M77").Value = -1839.5
$ws.Range("H116").Value = 506
$ws.Range("I116").Value = 356.88235
$ws.Range("J116").Value = 1013
$ws.Range("K116").Value = 356.88235
$ws.Range("L116").Value = 1013
$ws.Range("M116").Value = 1937.11765
$ws.Range("N116").Value = -5601
$ws.Range("H133").Value = 139995.8
$ws.Range("J133").Value = 139995.8
$ws.Range("L133").Value = 139995.8
$ws.Range("N133").Value = -145055.8

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 506
$ws.Range("I3").Value = 356.88235
$ws.Range("J3").Value = 1013
$ws.Range("K3").Value = 356.88235
$ws.Range("L3").Value = 1013
$ws.Range("M3").Value = -242.88235
$ws.Range("N3").Value = -1241

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 203633.8
$ws.Range("I31").Value = 252800
$ws.Range("K31").Value = 252800
$ws.Range("M31").Value = -252505
$ws.Range("H34").Value = 203633.8
$ws.Range("I34").Value = 252800
$ws.Range("K34").Value = 252800
$ws.Range("M34").Value = -252598
$ws.Range("H58").Value = 2035.138
$ws.Range("I58").Value = 1314.3914
$ws.Range("J58").Value = 4798
$ws.Range("K58").Value = 1314.3914
$ws.Range("L58").Value = 4798
$ws.Range("M58").Value = -1111.3914
$ws.Range("N58").Value = -5204
$ws.Range("H135").Value = 77294
$ws.Range("J135").Value = 77294
$ws.Range("L135").Value = 77294
$ws.Range("N135").Value = -87434
$ws.Range("H136").Value = 2035.138
$ws.Range("I136").Value = 1314.3914
$ws.Range("J136").Value = 4798
$ws.Range("K136").Value = 3943.1742
$ws.Range("L136").Value = 14394
$ws.Range("M136").Value = -1393.1742
$ws.Range("N136").Value = -19494

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 903.3333
$ws.Range("I34").Value = 130.66667
$ws.Range("J34").Value = 1160.8889
$ws.Range("K34").Value = 392.00001
$ws.Range("L34").Value = 3482.6667
$ws.Range("M34").Value = -308.00001
$ws.Range("N34").Value = -3650.6667
$ws.Range("H36").Value = 2408.1667
$ws.Range("I36").Value = 1290
$ws.Range("J36").Value = 7999
$ws.Range("K36").Value = 3870
$ws.Range("L36").Value = 23997
$ws.Range("M36").Value = -3701
$ws.Range("N36").Value = -24335
$ws.Range("H55").Value = 2098.5293
$ws.Range("I55").Value = 910
$ws.Range("J55").Value = 2464.2307
$ws.Range("K55").Value = 2730
$ws.Range("L55").Value = 7392.6921
$ws.Range("M55").Value = -2553
$ws.Range("N55").Value = -7746.6921
$ws.Range("H61").Value = 149.25
$ws.Range("I61").Value = 82.333336
$ws.Range("K61").Value = 247.000008
$ws.Range("M61").Value = -32.00000800000001
$ws.Range("H107").Value = 679.41174
$ws.Range("I107").Value = 472.77777
$ws.Range("J107").Value = 911.875
$ws.Range("K107").Value = 1418.33331
$ws.Range("L107").Value = 2735.625
$ws.Range("M107").Value = 501.66669
$ws.Range("N107").Value = -6575.625
$ws.Range("H113").Value = 3125825.5
$ws.Range("I113").Value = 4546188
$ws.Range("J113").Value = 1027.6
$ws.Range("K113").Value = 13638564
$ws.Range("L113").Value = 3082.8
$ws.Range("M113").Value = -13636394
$ws.Range("N113").Value = -7422.799999999999
$ws.Range("H122").Value = 1787.5238
$ws.Range("J122").Value = 1868.75
$ws.Range("L122").Value = 16818.75
$ws.Range("N122").Value = -21718.75

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 6518080
$ws.Range("I14").Value = 7482892
$ws.Range("K14").Value = 7482892
$ws.Range("M14").Value = -7482724
$ws.Range("H136").Value = 33356.926
$ws.Range("J136").Value = 33356.926
$ws.Range("L136").Value = 100070.778
$ws.Range("N136").Value = -105170.778

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6301.7144
$ws.Range("I7").Value = 5579.3335
$ws.Range("K7").Value = 5579.3335
$ws.Range("M7").Value = -5467.3335
$ws.Range("H35").Value = 530.4
$ws.Range("I35").Value = 530.4
$ws.Range("K35").Value = 530.4
$ws.Range("M35").Value = -194.4
$ws.Range("H126").Value = 6301.7144
$ws.Range("I126").Value = 5579.3335
$ws.Range("K126").Value = 16738.0005
$ws.Range("M126").Value = -14268.0005

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2192.4707
$ws.Range("I136").Value = 1330
$ws.Range("K136").Value = 3990
$ws.Range("M136").Value = -1440
